# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 currently carries the "bottom of table" border formatting (thicker
# bottom border etc.). Once rows 19-20 are removed below, row 18 becomes the
# new last row of the table, so copy that formatting onto row 18 first
# (while row 20 still exists) to preserve the closing border styling.
$ws.Range("B20:J20").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update summary figures at top of the account statement ---
# Valor Mora (total owed) recalculated for the updated set of periods
$ws.Range("E11").Value = 170820
# Cant. Periodos: now only 3 periods are being reported
$ws.Range("F13").Value = 3

# --- Update the 3 remaining period rows (16-18) ---
# New period labels, ascending order: 2504, 2505, 2506
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"

# Updated Salario Basico value for the reported periods
$ws.Range("G16").Value = 1450000
$ws.Range("G17").Value = 1450000
$ws.Range("G18").Value = 1450000

# --- Remove the two oldest period rows (formerly periods 2503 and 2502) ---
$ws.Rows("19:20").Delete()

Write-Output "done"
